# Auto-generated script applying crypto price/volume/coin updates
# Diff-driven update of Coin/Link/Price/Volume(1h) columns (rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.817.73'
$ws.Range("E2").Value = '  +1.25%  '
$ws.Range("D3").Value = '1.888.79'
$ws.Range("E3").Value = '  +1.76%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.48%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4758'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.63%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2895'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06590'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '99.50'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +18.17%  '
$ws.Range("D12").Value = '1.886.38'
$ws.Range("E12").Value = '  +1.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07597'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.144'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6628'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '309.32'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +34.76%  '
$ws.Range("D17").Value = '30.808.73'
$ws.Range("E17").Value = '  +1.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.22'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007593'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.64%  '
$ws.Range("D21").Value = '2.131.58'
$ws.Range("E21").Value = '  +2.04%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.127'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.235'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.327'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.48'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.54'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +14.70%  '
$ws.Range("E28").Value = '  +3.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1079'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.352'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.189'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.995'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05081'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.96%  '
$ws.Range("E34").Value = '  +2.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7310'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.717'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01962'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.705'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.079'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9108'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '108.18'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.000'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4220'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.634'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.17%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.11'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.39%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.399'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.65%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1231'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.10%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.017'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.83'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05637'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.394'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.97%  '

Write-Host "Applied all cell updates"
